$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.863.87"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.739.11"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.71"
$ws.Range("E5").Value = "  +4.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5205"
$ws.Range("E7").Value = "  -1.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2747"
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06158"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.744.43"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07184"
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.97"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6436"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.610"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.56"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.896.00"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.71"
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006764"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.966.36"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.286"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.638"
$ws.Range("E23").Value = "  -1.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.274"
$ws.Range("E24").Value = "  +2.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.74"
$ws.Range("E25").Value = "  -1.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.523"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.18"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.763"
$ws.Range("E28").Value = "  -1.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.11"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.918"
$ws.Range("E30").Value = "  +5.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08303"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.693"
$ws.Range("E32").Value = "  +4.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04626"
$ws.Range("E33").Value = "  +2.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.644"
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9891"
$ws.Range("E35").Value = "  +1.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6186"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.682"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01604"
$ws.Range("E38").Value = "  +1.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.927"
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9988"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.64"
$ws.Range("E41").Value = "  -2.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3848"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7425"
$ws.Range("E43").Value = "  +2.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.981"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1132"
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.244"
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05247"
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.79"
$ws.Range("E48").Value = "  +2.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.40"
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.601"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3413"
$ws.Range("E51").Value = "  +0.39%  "
